# Updates for 26 April
# Adds a new "4/25/20" column (AR) after the existing last column (AQ)
# on the single worksheet, copying formatting from column AQ and filling
# in the new day's death counts for each state row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day-of values for AR2:AR54, in row order (row 2 = Alabama ... row 54 = Wyoming)
$values = @(213,9,273,48,1695,672,1862,112,165,1055,907,5,14,56,1874,785,112,120,205,1703,50,875,2730,3274,244,221,283,14,53,204,60,5863,93,21908,305,16,711,194,87,1804,83,215,166,10,178,641,41,46,436,738,33,266,7)

# Header for the new column
$ws.Range("AR1").Value = " 4/25/20"

# Copy formatting (style) of the AQ column into the new AR column, then set values.
# Only rows 1-54 hold data in this column; rows 55/56 stay untouched (no AR cell there).
$ws.Range("AQ1:AQ54").Copy() | Out-Null
$ws.Range("AR1:AR54").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 44).Value = $values[$i]
}

# Update the active selection to match the new last column, as in the source workbook
$ws.Range("AR2").Select() | Out-Null
